$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row 72: Date, Hours, Task (new shared string "Bugfixes")
$ws.Range("A72").Value = 45391
$ws.Range("A72").NumberFormat = $ws.Range("A71").NumberFormat
$ws.Range("B72").Value = 5
$ws.Range("C72").Value = "Bugfixes"

# Update selection to match the diff
$ws.Range("C72").Select()
